$wb = $excel.ActiveWorkbook

# Update the "F" column (想去人数 / "want to go" count) on the "展览" sheet.
# Rows are keyed by their row number on that sheet.
$sheet1Updates = @{
    2  = 7566
    4  = 216
    5  = 12
    6  = 254
    7  = 1137
    8  = 198
    9  = 20
    10 = 139
    11 = 36
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

# Same events appear on "全部类型" too, but that sheet has one extra row
# (a concert entry, row 9) so the matching rows are shifted down by one
# from row 9 onward.
$sheet4Updates = @{
    2  = 7566
    4  = 216
    5  = 12
    6  = 254
    7  = 1137
    8  = 198
    10 = 20
    11 = 139
    12 = 36
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
